# "Admin control updates and customer enhancements"
# - Convert stored UTC timestamps to SGT (UTC+8) display strings across all
#   sheets (Users.Created, Products.Created, POS_Sales.Date,
#   Customer_Orders.Purchase_Date).
# - Decrement a few Products stock counts to reflect POS sales that already
#   happened (admin control updates).
# - Append new POS_Sales and Customer_Orders rows for orders placed since the
#   last export (customer enhancements).

$wb = $excel.ActiveWorkbook

# NOTE: this COM-interop host's [DateTime] cast / ParseExact do not produce a
# real System.DateTime (casting to [DateTime] is a no-op that leaves a
# string, and the resulting "string" has no .AddHours method) - so the
# UTC -> SGT (UTC+8) shift below is done with plain string/int arithmetic
# instead of DateTime APIs.
function Convert-ToSgt([string]$utcStamp) {
    $datePart = $utcStamp.Substring(0, 10)
    $timePart = $utcStamp.Substring(11, 8)

    $dateBits = $datePart.Split("-")
    $year = [int]$dateBits[0]
    $month = [int]$dateBits[1]
    $day = [int]$dateBits[2]

    $timeBits = $timePart.Split(":")
    $hour = [int]$timeBits[0]
    $minute = $timeBits[1]
    $second = $timeBits[2]

    $hour = $hour + 8
    if ($hour -ge 24) {
        $hour = $hour - 24
        $day = $day + 1

        $daysInMonth = 30
        if ($month -eq 1 -or $month -eq 3 -or $month -eq 5 -or $month -eq 7 -or $month -eq 8 -or $month -eq 10 -or $month -eq 12) {
            $daysInMonth = 31
        } elseif ($month -eq 2) {
            $daysInMonth = 28
            if (($year % 4 -eq 0 -and $year % 100 -ne 0) -or ($year % 400 -eq 0)) {
                $daysInMonth = 29
            }
        }

        if ($day -gt $daysInMonth) {
            $day = $day - $daysInMonth
            $month = $month + 1
            if ($month -gt 12) {
                $month = 1
                $year = $year + 1
            }
        }
    }

    $hourStr = $hour.ToString("D2")
    $monthStr = $month.ToString("D2")
    $dayStr = $day.ToString("D2")

    return "$year-$monthStr-$dayStr $hourStr`:$minute`:$second SGT"
}

# ---------------------------------------------------------------------
# Users sheet: shift "Created" timestamps to SGT
# ---------------------------------------------------------------------
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("E2").Value = Convert-ToSgt "2025-11-02 05:08:59"
$wsUsers.Range("E3").Value = Convert-ToSgt "2025-11-05 23:21:17"
$wsUsers.Range("E4").Value = Convert-ToSgt "2025-11-09 13:59:58"
$wsUsers.Range("E5").Value = Convert-ToSgt "2025-11-10 04:35:47"

# ---------------------------------------------------------------------
# Products sheet: shift "Created" timestamps to SGT + adjust stock counts
# ---------------------------------------------------------------------
$wsProducts = $wb.Worksheets.Item("Products")
$wsProducts.Range("H2").Value = Convert-ToSgt "2025-11-05 23:19:36"

$wsProducts.Range("E3").Value = 12
$wsProducts.Range("H3").Value = Convert-ToSgt "2025-11-08 06:47:45"

$wsProducts.Range("E4").Value = 17
$wsProducts.Range("H4").Value = Convert-ToSgt "2025-11-09 13:17:33"

$wsProducts.Range("E5").Value = 13
$wsProducts.Range("H5").Value = Convert-ToSgt "2025-11-09 13:17:53"

$wsProducts.Range("H6").Value = Convert-ToSgt "2025-11-09 13:33:38"
$wsProducts.Range("H7").Value = Convert-ToSgt "2025-11-09 13:34:32"
$wsProducts.Range("H8").Value = Convert-ToSgt "2025-11-09 13:49:55"

$wsProducts.Range("E9").Value = 8
$wsProducts.Range("H9").Value = Convert-ToSgt "2025-11-09 13:50:42"

$wsProducts.Range("H10").Value = Convert-ToSgt "2025-11-09 13:51:23"
$wsProducts.Range("H11").Value = Convert-ToSgt "2025-11-09 13:52:04"
$wsProducts.Range("H12").Value = Convert-ToSgt "2025-11-09 13:53:39"
$wsProducts.Range("H13").Value = Convert-ToSgt "2025-11-09 13:54:20"
$wsProducts.Range("H14").Value = Convert-ToSgt "2025-11-09 13:55:09"
$wsProducts.Range("H15").Value = Convert-ToSgt "2025-11-09 13:55:51"

# ---------------------------------------------------------------------
# POS_Sales sheet: shift "Date" timestamps to SGT + append new sales
# ---------------------------------------------------------------------
$wsPos = $wb.Worksheets.Item("POS_Sales")
$wsPos.Range("H2").Value = Convert-ToSgt "2025-11-05 23:19:59"
$wsPos.Range("H3").Value = Convert-ToSgt "2025-11-08 07:12:07"
$wsPos.Range("H4").Value = Convert-ToSgt "2025-11-08 07:25:20"
$wsPos.Range("H5").Value = Convert-ToSgt "2025-11-09 14:18:29"
$wsPos.Range("H6").Value = Convert-ToSgt "2025-11-09 15:37:44"
$wsPos.Range("H7").Value = Convert-ToSgt "2025-11-10 04:40:56"

$newPosSales = @(
    @{ Row=8;  A=7;  B=1; C=1399; D="cash"; E="None"; F=0; G='[{"product_id": 4, "name": "Eterno Saint", "price": 1399, "quantity": 1, "stock": 15}]';  H="2025-11-12 17:36:05 SGT" },
    @{ Row=9;  A=8;  B=1; C=1399; D="cash"; E="None"; F=0; G='[{"product_id": 4, "name": "Eterno Saint", "price": 1399, "quantity": 1, "stock": 15}]';  H="2025-11-12 17:36:05 SGT" },
    @{ Row=10; A=9;  B=1; C=1399; D="cash"; E="None"; F=0; G='[{"product_id": 4, "name": "Eterno Saint", "price": 1399, "quantity": 1, "stock": 15}]';  H="2025-11-12 17:36:05 SGT" },
    @{ Row=11; A=10; B=1; C=1399; D="cash"; E="None"; F=0; G='[{"product_id": 4, "name": "Eterno Saint", "price": 1399, "quantity": 1, "stock": 15}]';  H="2025-11-12 17:36:05 SGT" },
    @{ Row=12; A=11; B=1; C=1199; D="cash"; E="None"; F=0; G='[{"product_id": 3, "name": "Eterno Outlaw", "price": 1199, "quantity": 1, "stock": 18}]'; H="2025-11-12 17:47:37 SGT" }
)

foreach ($r in $newPosSales) {
    $wsPos.Cells.Item($r.Row, 1).Value = $r.A
    $wsPos.Cells.Item($r.Row, 2).Value = $r.B
    $wsPos.Cells.Item($r.Row, 3).Value = $r.C
    $wsPos.Cells.Item($r.Row, 4).Value = $r.D
    $wsPos.Cells.Item($r.Row, 5).Value = $r.E
    $wsPos.Cells.Item($r.Row, 6).Value = $r.F
    $wsPos.Cells.Item($r.Row, 7).Value = $r.G
    $wsPos.Cells.Item($r.Row, 8).Value = $r.H
}

# ---------------------------------------------------------------------
# Customer_Orders sheet: shift "Purchase_Date" timestamps to SGT + append
# new orders
# ---------------------------------------------------------------------
$wsOrders = $wb.Worksheets.Item("Customer_Orders")
$wsOrders.Range("L2").Value = Convert-ToSgt "2025-11-08 06:15:10"
$wsOrders.Range("L3").Value = Convert-ToSgt "2025-11-08 06:15:36"
$wsOrders.Range("L4").Value = Convert-ToSgt "2025-11-08 07:10:04"
$wsOrders.Range("L5").Value = Convert-ToSgt "2025-11-09 12:09:55"
$wsOrders.Range("L6").Value = Convert-ToSgt "2025-11-09 12:12:14"
$wsOrders.Range("L7").Value = Convert-ToSgt "2025-11-09 12:12:59"
$wsOrders.Range("L8").Value = Convert-ToSgt "2025-11-09 12:13:27"
$wsOrders.Range("L9").Value = Convert-ToSgt "2025-11-09 14:01:41"
$wsOrders.Range("L10").Value = Convert-ToSgt "2025-11-09 14:16:10"
$wsOrders.Range("L11").Value = Convert-ToSgt "2025-11-09 15:36:22"
$wsOrders.Range("L12").Value = Convert-ToSgt "2025-11-10 02:14:24"
$wsOrders.Range("L13").Value = Convert-ToSgt "2025-11-10 04:37:35"

$address15 = "Address 1:`nBuilding Number: 56`nStreet Name: Rue Drummond`nStreet Address: Ritz-Carlton Montreal`nState: Quebec`nCity: Montreal`nPost Code: H3G 1Y9"

$newOrders = @(
    @{ Row=14; A=13; B=3; C="kaizen"; D="boarratjabol@gmail.com"; E="123123123";  F=799;  G=65; H=864;  I="credit_card"; J="completed"; K='[{"product_id": 2, "product_name": "Eterno Grace", "quantity": 1, "price": 799.0}]'; L="2025-11-12 17:38:40 SGT" },
    @{ Row=15; A=14; B=3; C="kaizen"; D="boarratjabol@gmail.com"; E=$address15; F=4599; G=61; H=4660; I="gcash";       J="completed"; K='[{"product_id": 8, "product_name": "Eterno Abyss", "quantity": 1, "price": 4599.0}]'; L="2025-11-12 17:49:27 SGT" }
)

foreach ($o in $newOrders) {
    $wsOrders.Cells.Item($o.Row, 1).Value = $o.A
    $wsOrders.Cells.Item($o.Row, 2).Value = $o.B
    $wsOrders.Cells.Item($o.Row, 3).Value = $o.C
    $wsOrders.Cells.Item($o.Row, 4).Value = $o.D
    $wsOrders.Cells.Item($o.Row, 5).Value = $o.E
    $wsOrders.Cells.Item($o.Row, 6).Value = $o.F
    $wsOrders.Cells.Item($o.Row, 7).Value = $o.G
    $wsOrders.Cells.Item($o.Row, 8).Value = $o.H
    $wsOrders.Cells.Item($o.Row, 9).Value = $o.I
    $wsOrders.Cells.Item($o.Row, 10).Value = $o.J
    $wsOrders.Cells.Item($o.Row, 11).Value = $o.K
    $wsOrders.Cells.Item($o.Row, 12).Value = $o.L
}
